$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.088.93"
$ws.Range("E2").Value = "  -1.96%  "

# Row 3
$ws.Range("D3").Value = "1.797.78"
$ws.Range("E3").Value = "  -2.81%  "

# Row 4
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.47%  "

# Row 5
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'308.13"
$ws.Range("E5").Value = "  -1.75%  "

# Row 6
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  +0.69%  "

# Row 7
$ws.Range("D7").Value = "'0.4174"
$ws.Range("E7").Value = "  -2.31%  "

# Row 8
$ws.Range("D8").Value = "'0.3549"
$ws.Range("E8").Value = "  -3.89%  "

# Row 9
$ws.Range("D9").Value = "'0.07063"
$ws.Range("E9").Value = "  -3.49%  "

# Row 10
$ws.Range("D10").Value = "'0.8426"
$ws.Range("E10").Value = "  -3.91%  "

# Row 11
$ws.Range("D11").Value = "'20.18"
$ws.Range("E11").Value = "  -3.64%  "

# Row 12
$ws.Range("D12").Value = "1.719.41"
$ws.Range("E12").Value = "  -7.42%  "

# Row 13
$ws.Range("D13").Value = "'5.271"
$ws.Range("E13").Value = "  -2.94%  "

# Row 14
$ws.Range("D14").Value = "'6.335"
$ws.Range("E14").Value = "  -3.48%  "

# Row 15
$ws.Range("D15").Value = "'0.06759"
$ws.Range("E15").Value = "  -2.63%  "

# Row 16
$ws.Range("E16").Value = "  +0.73%  "

# Row 17
$ws.Range("D17").Value = "'79.90"
$ws.Range("E17").Value = "  -0.79%  "

# Row 18
$ws.Range("D18").Value = "'0.000008698"
$ws.Range("E18").Value = "  -4.01%  "

# Row 19
$ws.Range("D19").Value = "'1.005"
$ws.Range("E19").Value = "  +0.79%  "

# Row 20
$ws.Range("D20").Value = "'15.01"
$ws.Range("E20").Value = "  -3.29%  "

# Row 21
$ws.Range("D21").Value = "27.184.67"
$ws.Range("E21").Value = "  -1.81%  "

# Row 22
$ws.Range("D22").Value = "'5.052"
$ws.Range("E22").Value = "  -0.74%  "

# Row 23
$ws.Range("D23").Value = "'10.91"
$ws.Range("E23").Value = "  -0.02%  "

# Row 24
$ws.Range("D24").Value = "2.036.93"
$ws.Range("E24").Value = "  -2.10%  "

# Row 25
$ws.Range("D25").Value = "'1.954"
$ws.Range("E25").Value = "  -0.22%  "

# Row 26
$ws.Range("D26").Value = "'153.16"

# Row 27
$ws.Range("E27").Value = "  -2.21%  "

# Row 28
$ws.Range("D28").Value = "'5.004"
$ws.Range("E28").Value = "  -5.13%  "

# Row 29
$ws.Range("D29").Value = "'112.59"
$ws.Range("E29").Value = "  -2.44%  "

# Row 30
$ws.Range("D30").Value = "'1.645"
$ws.Range("E30").Value = "  -11.93%  "

# Row 31
$ws.Range("D31").Value = "'0.08890"
$ws.Range("E31").Value = "  -0.27%  "

# Row 32
$ws.Range("D32").Value = "'0.7178"
$ws.Range("E32").Value = "  -8.40%  "

# Row 33
$ws.Range("D33").Value = "'2.862"
$ws.Range("E33").Value = "  -4.06%  "

# Row 34
$ws.Range("D34").Value = "'4.321"
$ws.Range("E34").Value = "  -5.84%  "

# Row 35
$ws.Range("D35").Value = "'1.005"
$ws.Range("E35").Value = "  +0.73%  "

# Row 36
$ws.Range("E36").Value = "  -7.17%  "

# Row 37
$ws.Range("D37").Value = "'1.072"
$ws.Range("E37").Value = "  -2.96%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.05105"
$ws.Range("E38").Value = "  -5.75%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01891"
$ws.Range("E39").Value = "  -3.14%  "

# Row 40
$ws.Range("D40").Value = "'0.1614"
$ws.Range("E40").Value = "  -3.46%  "

# Row 41
$ws.Range("D41").Value = "'0.4928"
$ws.Range("E41").Value = "  -4.26%  "

# Row 42
$ws.Range("D42").Value = "'2.583"
$ws.Range("E42").Value = "  -8.95%  "

# Row 43
$ws.Range("D43").Value = "'6.108"
$ws.Range("E43").Value = "  -9.31%  "

# Row 44
$ws.Range("D44").Value = "'8.033"
$ws.Range("E44").Value = "  -7.41%  "

# Row 45
$ws.Range("D45").Value = "'1.005"
$ws.Range("E45").Value = "  +0.90%  "

# Row 46
$ws.Range("D46").Value = "'104.23"
$ws.Range("E46").Value = "  -2.86%  "

# Row 47
$ws.Range("D47").Value = "'10.17"
$ws.Range("E47").Value = "  -3.53%  "

# Row 48
$ws.Range("D48").Value = "'0.06308"
$ws.Range("E48").Value = "  -3.79%  "

# Row 49
$ws.Range("D49").Value = "'0.4502"
$ws.Range("E49").Value = "  -5.17%  "

# Row 50
$ws.Range("D50").Value = "'1.588"
$ws.Range("E50").Value = "  -3.73%  "

# Row 51
$ws.Range("D51").Value = "'61.77"
$ws.Range("E51").Value = "  -4.55%  "
